$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the "desc" value for the soldier card row (F4), matching the new
# shared string "SOLDIER_DESC_1" added for this row's description key.
$ws.Range("F4").Value = "SOLDIER_DESC_1"

# Move/record the active selection as it ended up after the edit.
$ws.Range("H7").Select()
